$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.324.11'
$ws.Range('D3').Value = '1.667.30'
$ws.Range('E3').Value = '  +0.84%  '
$ws.Range('D4').Value = '1.009'
$ws.Range('E4').Value = '  +0.35%  '
$ws.Range('D5').Value = '220.74'
$ws.Range('E5').Value = '  +1.43%  '
$ws.Range('D6').Value = '0.5309'
$ws.Range('E6').Value = '  +0.16%  '
$ws.Range('E7').Value = '  +0.33%  '
$ws.Range('D8').Value = '0.2649'
$ws.Range('E8').Value = '  +1.11%  '
$ws.Range('D9').Value = '0.06365'
$ws.Range('E9').Value = '  +0.53%  '
$ws.Range('D10').Value = '20.86'
$ws.Range('E10').Value = '  +2.36%  '
$ws.Range('D11').Value = '0.07842'
$ws.Range('E11').Value = '  +0.51%  '
$ws.Range('D12').Value = '4.516'
$ws.Range('E12').Value = '  -0.01%  '
$ws.Range('D13').Value = '1.669.34'
$ws.Range('E13').Value = '  +2.33%  '
$ws.Range('D14').Value = '1.895.82'
$ws.Range('E14').Value = '  +0.78%  '
$ws.Range('D15').Value = '0.5596'
$ws.Range('E15').Value = '  +1.92%  '
$ws.Range('D16').Value = '0.0₅8161'
$ws.Range('E16').Value = '  -0.19%  '
$ws.Range('D17').Value = '65.77'
$ws.Range('E17').Value = '  +0.46%  '
$ws.Range('D18').Value = '26.328.00'
$ws.Range('E18').Value = '  +0.73%  '
$ws.Range('E19').Value = '  +0.41%  '
$ws.Range('E20').Value = '  +2.55%  '
$ws.Range('D21').Value = '196.61'
$ws.Range('E21').Value = '  +2.90%  '
$ws.Range('D23').Value = '6.044'
$ws.Range('E23').Value = '  +0.67%  '
$ws.Range('D24').Value = '1.010'
$ws.Range('E24').Value = '  +0.31%  '
$ws.Range('D25').Value = '145.66'
$ws.Range('E25').Value = '  +0.29%  '
$ws.Range('E26').Value = '  -0.34%  '
$ws.Range('D27').Value = '7.238'
$ws.Range('E27').Value = '  +0.56%  '
$ws.Range('E28').Value = '  +1.06%  '
$ws.Range('D29').Value = '1.508'
$ws.Range('E29').Value = '  +2.17%  '
$ws.Range('D30').Value = '0.05885'
$ws.Range('E30').Value = '  +2.76%  '
$ws.Range('D31').Value = '1.286'
$ws.Range('E31').Value = '  +1.03%  '
$ws.Range('E32').Value = '  -0.28%  '
$ws.Range('D33').Value = '3.333'
$ws.Range('E33').Value = '  +2.19%  '
$ws.Range('E34').Value = '  +1.14%  '
$ws.Range('D35').Value = '2.831'
$ws.Range('E35').Value = '  +0.84%  '
$ws.Range('D36').Value = '0.9597'
$ws.Range('E36').Value = '  +1.23%  '
$ws.Range('E37').Value = '  +0.56%  '
$ws.Range('E38').Value = '  +1.31%  '
$ws.Range('D39').Value = '0.01613'
$ws.Range('E39').Value = '  +0.81%  '
$ws.Range('D40').Value = '5.944'
$ws.Range('E40').Value = '  +2.89%  '
$ws.Range('D41').Value = '1.074.88'
$ws.Range('E41').Value = '  +3.42%  '
$ws.Range('D42').Value = '0.8620'
$ws.Range('E42').Value = '  +1.58%  '
$ws.Range('E43').Value = '  +0.36%  '
$ws.Range('D44').Value = '102.71'
$ws.Range('E44').Value = '  -1.12%  '
$ws.Range('D45').Value = '1.806.03'
$ws.Range('E45').Value = '  +0.70%  '
$ws.Range('D46').Value = '58.38'
$ws.Range('E46').Value = '  +2.71%  '
$ws.Range('D47').Value = '0.0₈106'
$ws.Range('E47').Value = '  +1.25%  '
$ws.Range('D48').Value = '1.015'
$ws.Range('E48').Value = '  +0.87%  '
$ws.Range('D49').Value = '0.4411'
$ws.Range('E49').Value = '  +1.25%  '
$ws.Range('D50').Value = '8.033'
$ws.Range('E50').Value = '  +2.05%  '
$ws.Range('E51').Value = '  -0.08%  '
